# Weekly update: insert a new week of price data (2 new rows) at the top
# of the existing "Pepino ensalada" data block (rows 973-1044), pushing the
# previously-existing rows down by two positions (-> 975-1046).
#
# New rows (date 2023-12-05, serial 45265):
#   row 973: Calidad "Primera"
#   row 974: Calidad "Segunda"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 973, shifting the
# existing data (973-1044) down to (975-1046).
$ws.Rows.Item(973).Resize(2).Insert()

# ---- New row 973 ("Primera") ----
$ws.Cells.Item(973, 1).Value = 8
$ws.Cells.Item(973, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(973, 3).Value = "Coquimbo"
$ws.Cells.Item(973, 4).Value = 45265
$ws.Cells.Item(973, 5).Value = 4
$ws.Cells.Item(973, 6).Value = 100112043
$ws.Cells.Item(973, 7).Value = "Pepino ensalada"
$ws.Cells.Item(973, 8).Value = "Sin especificar"
$ws.Cells.Item(973, 9).Value = "Primera"
$ws.Cells.Item(973, 10).Value = 600
$ws.Cells.Item(973, 11).Value = 14000
$ws.Cells.Item(973, 12).Value = 15000
$ws.Cells.Item(973, 13).Value = 14500
$ws.Cells.Item(973, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(973, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(973, 16).Value = 242
$ws.Cells.Item(973, 17).Value = 60
$ws.Cells.Item(973, 18).Value = "Hortaliza"

# ---- New row 974 ("Segunda") ----
$ws.Cells.Item(974, 1).Value = 8
$ws.Cells.Item(974, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(974, 3).Value = "Coquimbo"
$ws.Cells.Item(974, 4).Value = 45265
$ws.Cells.Item(974, 5).Value = 4
$ws.Cells.Item(974, 6).Value = 100112043
$ws.Cells.Item(974, 7).Value = "Pepino ensalada"
$ws.Cells.Item(974, 8).Value = "Sin especificar"
$ws.Cells.Item(974, 9).Value = "Segunda"
$ws.Cells.Item(974, 10).Value = 400
$ws.Cells.Item(974, 11).Value = 10000
$ws.Cells.Item(974, 12).Value = 11000
$ws.Cells.Item(974, 13).Value = 10500
$ws.Cells.Item(974, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(974, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(974, 16).Value = 131
$ws.Cells.Item(974, 17).Value = 80
$ws.Cells.Item(974, 18).Value = "Hortaliza"
